# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" quarter sheet right after the "总计" (summary) sheet,
#    pushing the existing quarter sheets ("2022-Q3" .. "2021-Q3") down by one
#    position. The new sheet is created by copying the existing "2022-Q3"
#    sheet (2nd sheet) as a template so it inherits the exact same layout,
#    styles, page setup, etc., then trimming it down to the two new fund rows
#    and overwriting their values.
# 2) Update the "总计" (summary) sheet: add a new data row for "2022-Q4"
#    right under the header, shifting the previous quarter rows down by one,
#    and append the row that drops off the end ("2021-Q3") as the new last row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: new "2022-Q4" worksheet with fund holdings data
# ---------------------------------------------------------------------------

$template = $wb.Worksheets.Item(2)      # currently "2022-Q3"
$template.Copy($template)               # places the copy right before itself

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template sheet had 4 data rows (rows 2-5); we only need 2, so drop the
# extra two rows (rows 4-5) to get dimension back down to A1:H3.
$newSheet.Rows.Item(4).Resize(2).Delete()

# Columns B and D:G hold numeric-looking data that must stay TEXT (to keep
# leading zeros on fund codes and trailing zeros in the percentages), so mark
# the range as text before writing, then restore the style afterwards so we
# don't leave a stray number-format style behind.
$textRange = $newSheet.Range("B2:G3")
$textRange.NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "013166"
$newSheet.Range("C2").Value = "东兴宸祥量化混合A"
$newSheet.Range("D2").Value = "0.38"
$newSheet.Range("E2").Value = "93.88"
$newSheet.Range("F2").Value = "1.06"
$newSheet.Range("G2").Value = "0.0040"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "013167"
$newSheet.Range("C3").Value = "东兴宸祥量化混合C"
$newSheet.Range("D3").Value = "0.08"
$newSheet.Range("E3").Value = "93.88"
$newSheet.Range("F3").Value = "1.06"
$newSheet.Range("G3").Value = "0.0008"
$newSheet.Range("H3").Value = 8

$textRange.Style = "Normal"

# ---------------------------------------------------------------------------
# Part 2: update the "总计" summary sheet with the new quarter's row
# ---------------------------------------------------------------------------

$summary = $wb.Worksheets.Item(1)       # "总计"

$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 4).Value = 0

# Fix up styling: the inserted row pulled its look from the header row, so
# re-apply the normal bordered/centred "index" style to A2 from the row below
# it, and clear the header-ish style that leaked onto B2:D2.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)     # xlPasteFormats

# The "index" column (A) is just the sequential 0-based row counter; restore
# it for every row now that a new row 2 was inserted.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5

$excel.CutCopyMode = 0
